# Generate Report for Handoff
# Adds one new tracked file ("58631033-...") to each of the three tables
# (Overview, zh-cn, de-de) as a new row, mirroring the existing
# "b5d706ad-..." entry, wires up its hyperlink, and widens a couple of
# date/time columns to fit the new timestamp text.

$wb = $excel.ActiveWorkbook

$newMd        = "58631033-4500-40ee-a97a-c4864d134564ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdPath    = "e2e\58631033-4500-40ee-a97a-c4864d134564ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9743a36d3df3136c7b3d4f0dcda9ba3ee0af8978/e2e/58631033-4500-40ee-a97a-c4864d134564ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

$zhXlf        = "58631033-4500-40ee-a97a-c4864d134564oooooooooooooooooooooooooooooooooooooooo.140cfda764fe0a1474d9073622b74c921d680328.zh-cn.xlf"
$deXlf        = "58631033-4500-40ee-a97a-c4864d134564oooooooooooooooooooooooooooooooooooooooo.140cfda764fe0a1474d9073622b74c921d680328.de-de.xlf"

$dateFmt      = "yyyy-mm-dd HH:mm:ss"
$handoffDate  = "2016-08-23 22:25:29"
$zhXlfDate    = "2016-08-23 22:25:24"
$noDate       = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet: one summary row per tracked file.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$rowOv = $loOverview.ListRows.Count + 1

$wsOverview.Cells.Item($rowOv, 1).Value = $newMd
$wsOverview.Cells.Item($rowOv, 2).Value = $newMdPath
$wsOverview.Cells.Item($rowOv, 3).Value = ".md"
$wsOverview.Cells.Item($rowOv, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($rowOv, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($rowOv, 7).Value = $handoffDate
$wsOverview.Cells.Item($rowOv, 7).NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B" + $rowOv), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMdPath) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: same file, zh-cn handoff xliff.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$rowZh = $loZh.ListRows.Count + 1

$wsZh.Cells.Item($rowZh, 1).Value  = $newMd          # Source File Name
$wsZh.Cells.Item($rowZh, 2).Value  = ".md"            # File Extension
$wsZh.Cells.Item($rowZh, 3).Value  = "Ready for handoff" # Status
$wsZh.Cells.Item($rowZh, 4).Value  = "e2e"            # Source Path
$wsZh.Cells.Item($rowZh, 5).Value  = "ht"             # Priority
$wsZh.Cells.Item($rowZh, 6).Value  = "False"          # Content Duplicate
$wsZh.Cells.Item($rowZh, 7).Value  = $zhXlf           # Latest Handoff File
$wsZh.Cells.Item($rowZh, 8).Value  = $zhXlfDate       # Latest Handoff Datetime
$wsZh.Cells.Item($rowZh, 8).NumberFormat = $dateFmt
$wsZh.Cells.Item($rowZh, 9).Value  = ""               # Latest Target File
$wsZh.Cells.Item($rowZh, 10).Value = ""               # Latest Handback File
$wsZh.Cells.Item($rowZh, 11).Value = $noDate          # Latest Handback DateTime
$wsZh.Cells.Item($rowZh, 11).NumberFormat = $dateFmt
$wsZh.Cells.Item($rowZh, 12).Value = ""               # Reference Tokens
$wsZh.Cells.Item($rowZh, 13).Value = "True"           # To be localized
$wsZh.Cells.Item($rowZh, 14).Value = ""               # Dependency From
$wsZh.Cells.Item($rowZh, 15).Value = "False"          # Has metadata
$wsZh.Cells.Item($rowZh, 16).Value = ""               # Error Detail

$wsZh.Hyperlinks.Add($wsZh.Range("A" + $rowZh), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: same file, de-de handoff xliff.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$rowDe = $loDe.ListRows.Count + 1

$wsDe.Cells.Item($rowDe, 1).Value  = $newMd          # Source File Name
$wsDe.Cells.Item($rowDe, 2).Value  = ".md"            # File Extension
$wsDe.Cells.Item($rowDe, 3).Value  = "Ready for handoff" # Status
$wsDe.Cells.Item($rowDe, 4).Value  = "e2e"            # Source Path
$wsDe.Cells.Item($rowDe, 5).Value  = "ht"             # Priority
$wsDe.Cells.Item($rowDe, 6).Value  = "False"          # Content Duplicate
$wsDe.Cells.Item($rowDe, 7).Value  = $deXlf           # Latest Handoff File
$wsDe.Cells.Item($rowDe, 8).Value  = $handoffDate     # Latest Handoff Datetime
$wsDe.Cells.Item($rowDe, 8).NumberFormat = $dateFmt
$wsDe.Cells.Item($rowDe, 9).Value  = ""               # Latest Target File
$wsDe.Cells.Item($rowDe, 10).Value = ""               # Latest Handback File
$wsDe.Cells.Item($rowDe, 11).Value = $noDate          # Latest Handback DateTime
$wsDe.Cells.Item($rowDe, 11).NumberFormat = $dateFmt
$wsDe.Cells.Item($rowDe, 12).Value = ""               # Reference Tokens
$wsDe.Cells.Item($rowDe, 13).Value = "True"           # To be localized
$wsDe.Cells.Item($rowDe, 14).Value = ""               # Dependency From
$wsDe.Cells.Item($rowDe, 15).Value = "False"          # Has metadata
$wsDe.Cells.Item($rowDe, 16).Value = ""               # Error Detail

$wsDe.Hyperlinks.Add($wsDe.Range("A" + $rowDe), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null

# ---------------------------------------------------------------------
# Widen the "Latest Handoff Datetime" / Latest Handback DateTime style
# columns so the longer timestamps introduced above are fully visible
# (matches column width bump seen for these date columns). The engine
# quantizes ColumnWidth to whole pixels, so 16.33 is the closest input
# that lands on the target stored width.
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZh.Columns.Item(3).ColumnWidth = 16.33
$wsDe.Columns.Item(3).ColumnWidth = 16.33

Write-Host "Handoff report row added to Overview, zh-cn, and de-de sheets."
